$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1373
$ws.Range("I19").Value = 500
$ws.Range("J19").Value = 1664
$ws.Range("K19").Value = 500
$ws.Range("L19").Value = 1664
$ws.Range("M19").Value = -325
$ws.Range("N19").Value = -2014
$ws.Range("H62").Value = 13890442
$ws.Range("I62").Value = 18520014
$ws.Range("K62").Value = 18520014
$ws.Range("M62").Value = -18519390
$ws.Range("H65").Value = 13890442
$ws.Range("I65").Value = 18520014
$ws.Range("K65").Value = 92600070
$ws.Range("M65").Value = -92596950
$ws.Range("H98").Value = 1225.3549
$ws.Range("I98").Value = 951.6
$ws.Range("J98").Value = 2366
$ws.Range("K98").Value = 951.6
$ws.Range("L98").Value = 2366
$ws.Range("M98").Value = 546.4
$ws.Range("N98").Value = -5362
$ws.Range("H122").Value = 1225.3549
$ws.Range("I122").Value = 951.6
$ws.Range("J122").Value = 2366
$ws.Range("K122").Value = 2854.8
$ws.Range("L122").Value = 7098
$ws.Range("M122").Value = -404.8000000000002
$ws.Range("N122").Value = -11998
$ws.Range("H137").Value = 1281.1091
$ws.Range("J137").Value = 1787.5294
$ws.Range("L137").Value = 5362.5882
$ws.Range("N137").Value = -10462.5882

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 840.7222
$ws.Range("I2").Value = 860
$ws.Range("J2").Value = 513
$ws.Range("K2").Value = 860
$ws.Range("L2").Value = 513
$ws.Range("M2").Value = -747
$ws.Range("N2").Value = -739
$ws.Range("H32").Value = 4956.64
$ws.Range("I32").Value = 3396.6155
$ws.Range("J32").Value = 10487.637
$ws.Range("K32").Value = 3396.6155
$ws.Range("L32").Value = 10487.637
$ws.Range("M32").Value = -3109.6155
$ws.Range("N32").Value = -11061.637
$ws.Range("H45").Value = 10171.272
$ws.Range("I45").Value = 10171.272
$ws.Range("K45").Value = 10171.272
$ws.Range("M45").Value = -9794.272000000001
$ws.Range("H116").Value = 840.7222
$ws.Range("I116").Value = 860
$ws.Range("J116").Value = 513
$ws.Range("K116").Value = 860
$ws.Range("L116").Value = 513
$ws.Range("M116").Value = 1434
$ws.Range("N116").Value = -5101
$ws.Range("H122").Value = 1832392.2
$ws.Range("I122").Value = 2331708.8
$ws.Range("J122").Value = 1565
$ws.Range("K122").Value = 6995126.399999999
$ws.Range("L122").Value = 4695
$ws.Range("M122").Value = -6992676.399999999
$ws.Range("N122").Value = -9595
$ws.Range("H132").Value = 2130408.5
$ws.Range("I132").Value = 1890.8182
$ws.Range("J132").Value = 7147629
$ws.Range("K132").Value = 5672.4546
$ws.Range("L132").Value = 21442887
$ws.Range("M132").Value = -3142.4546
$ws.Range("N132").Value = -21447947

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 840.7222
$ws.Range("I3").Value = 860
$ws.Range("J3").Value = 513
$ws.Range("K3").Value = 860
$ws.Range("L3").Value = 513
$ws.Range("M3").Value = -746
$ws.Range("N3").Value = -741
$ws.Range("H105").Value = 14293194
$ws.Range("I105").Value = 29424910
$ws.Range("J105").Value = 2128.889
$ws.Range("K105").Value = 29424910
$ws.Range("L105").Value = 2128.889
$ws.Range("M105").Value = -29423163
$ws.Range("N105").Value = -5622.889
$ws.Range("H134").Value = 5250.943
$ws.Range("I134").Value = 7932.7646
$ws.Range("J134").Value = 2718.111
$ws.Range("K134").Value = 23798.2938
$ws.Range("L134").Value = 8154.333
$ws.Range("M134").Value = -21263.2938
$ws.Range("N134").Value = -13224.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1625.8125
$ws.Range("I16").Value = 1548.75
$ws.Range("J16").Value = 1702.875
$ws.Range("K16").Value = 1548.75
$ws.Range("L16").Value = 1702.875
$ws.Range("M16").Value = -1261.75
$ws.Range("N16").Value = -2276.875
$ws.Range("H22").Value = 352.625
$ws.Range("J22").Value = 458.2
$ws.Range("L22").Value = 458.2
$ws.Range("N22").Value = -1158.2
$ws.Range("H31").Value = 302626.44
$ws.Range("I31").Value = 1818.4517
$ws.Range("J31").Value = 746676.3
$ws.Range("K31").Value = 1818.4517
$ws.Range("L31").Value = 746676.3
$ws.Range("M31").Value = -1523.4517
$ws.Range("N31").Value = -747266.3
$ws.Range("H34").Value = 302626.44
$ws.Range("I34").Value = 1818.4517
$ws.Range("J34").Value = 746676.3
$ws.Range("K34").Value = 1818.4517
$ws.Range("L34").Value = 746676.3
$ws.Range("M34").Value = -1616.4517
$ws.Range("N34").Value = -747080.3
$ws.Range("H113").Value = 1625.8125
$ws.Range("I113").Value = 1548.75
$ws.Range("J113").Value = 1702.875
$ws.Range("K113").Value = 1548.75
$ws.Range("L113").Value = 1702.875
$ws.Range("M113").Value = 621.25
$ws.Range("N113").Value = -6042.875
$ws.Range("H122").Value = 1677
$ws.Range("I122").Value = 1215.5
$ws.Range("J122").Value = 2138.5
$ws.Range("K122").Value = 3646.5
$ws.Range("L122").Value = 6415.5
$ws.Range("M122").Value = -1196.5
$ws.Range("N122").Value = -11315.5
$ws.Range("H132").Value = 2274.608
$ws.Range("I132").Value = 1931.3235
$ws.Range("J132").Value = 2961.1765
$ws.Range("K132").Value = 5793.970499999999
$ws.Range("L132").Value = 8883.529500000001
$ws.Range("M132").Value = -3263.970499999999
$ws.Range("N132").Value = -13943.5295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 20000.666
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("I24").Value = 10751975
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 10751975
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -10751802
$ws.Range("N24").ClearContents()
$ws.Range("H30").Value = 20000.666
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H93").Value = 9251
$ws.Range("J93").Value = 9251
$ws.Range("L93").Value = 9251
$ws.Range("N93").Value = -12995
$ws.Range("H102").Value = 1375.0714
$ws.Range("I102").Value = 1177.5
$ws.Range("J102").Value = 1730.7
$ws.Range("K102").Value = 1177.5
$ws.Range("L102").Value = 1730.7
$ws.Range("M102").Value = 444.5
$ws.Range("N102").Value = -4974.7
$ws.Range("I113").Value = 83334136
$ws.Range("J113").Value = 1336.6666
$ws.Range("K113").Value = 83334136
$ws.Range("L113").Value = 1336.6666
$ws.Range("M113").Value = -83331966
$ws.Range("N113").Value = -5676.6666
$ws.Range("H122").Value = 42594544
$ws.Range("I122").Value = 70989280
$ws.Range("J122").Value = 2442.3
$ws.Range("K122").Value = 212967840
$ws.Range("L122").Value = 7326.900000000001
$ws.Range("M122").Value = -212965390
$ws.Range("N122").Value = -12226.9
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 2379.7837
$ws.Range("I132").Value = 1814.1818
$ws.Range("J132").Value = 3209.3333
$ws.Range("K132").Value = 5442.5454
$ws.Range("L132").Value = 9627.999899999999
$ws.Range("M132").Value = -2912.5454
$ws.Range("N132").Value = -14687.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 168900.67
$ws.Range("I7").Value = 201780.8
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 201780.8
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -201668.8
$ws.Range("N7").Value = -4724
$ws.Range("H40").Value = 142863420
$ws.Range("I40").Value = 250002260
$ws.Range("J40").Value = 11666.667
$ws.Range("K40").Value = 250002260
$ws.Range("L40").Value = 11666.667
$ws.Range("M40").Value = -250002124
$ws.Range("N40").Value = -11938.667
$ws.Range("H46").Value = 20834156
$ws.Range("J46").Value = 1060.625
$ws.Range("L46").Value = 1060.625
$ws.Range("N46").Value = -1436.625
$ws.Range("H55").Value = 150000300
$ws.Range("I55").Value = 200000260
$ws.Range("J55").Value = 100000360
$ws.Range("K55").Value = 200000260
$ws.Range("L55").Value = 100000360
$ws.Range("M55").Value = -200000087
$ws.Range("N55").Value = -100000706
$ws.Range("H82").Value = 1003498.2
$ws.Range("I82").Value = 1667633.4
$ws.Range("J82").Value = 206536
$ws.Range("K82").Value = 1667633.4
$ws.Range("L82").Value = 206536
$ws.Range("M82").Value = -1667272.4
$ws.Range("N82").Value = -207258
$ws.Range("H85").Value = 1003498.2
$ws.Range("I85").Value = 1667633.4
$ws.Range("J85").Value = 206536
$ws.Range("K85").Value = 1667633.4
$ws.Range("L85").Value = 206536
$ws.Range("M85").Value = -1666385.4
$ws.Range("N85").Value = -209032
$ws.Range("H122").Value = 3136756.8
$ws.Range("I122").Value = 3765014.5
$ws.Range("J122").Value = 1431485.8
$ws.Range("K122").Value = 11295043.5
$ws.Range("L122").Value = 4294457.4
$ws.Range("M122").Value = -11292593.5
$ws.Range("N122").Value = -4299357.4
$ws.Range("H126").Value = 168900.67
$ws.Range("I126").Value = 201780.8
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 605342.3999999999
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -602872.3999999999
$ws.Range("N126").Value = -18440
$ws.Range("H132").Value = 8763053
$ws.Range("I132").Value = 11372396
$ws.Range("J132").Value = 3114.8572
$ws.Range("K132").Value = 34117188
$ws.Range("L132").Value = 9344.571599999999
$ws.Range("M132").Value = -34114658
$ws.Range("N132").Value = -14404.5716
$ws.Range("H136").Value = 7062.3335
$ws.Range("I136").Value = 4377.857
$ws.Range("K136").Value = 13133.571
$ws.Range("M136").Value = -10583.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2236.5
$ws.Range("I96").Value = 2600
$ws.Range("J96").Value = 1873
$ws.Range("K96").Value = 2600
$ws.Range("L96").Value = 1873
$ws.Range("M96").Value = -1227
$ws.Range("N96").Value = -4619
$ws.Range("H122").Value = 1885.3429
$ws.Range("I122").Value = 1857.1923
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 5571.5769
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -3121.5769
$ws.Range("N122").Value = -10799.9998
$ws.Range("H126").Value = 947.3333
$ws.Range("I126").Value = 636.8
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 1910.4
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = 559.6000000000001
$ws.Range("N126").Value = -12440
$ws.Range("H136").Value = 7939051.5
$ws.Range("I136").Value = 2747.2307
$ws.Range("K136").Value = 8241.6921
$ws.Range("M136").Value = -5691.6921

Write-Host "Applied all changes"